$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.236.32'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7095'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +0.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '237.80'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -0.57%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.0000'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.08187'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  +10.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3039'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -0.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.28'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -0.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08191'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.845.74'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  -0.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.176'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7087'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -2.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.52'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +0.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.218.36'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007922'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +3.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.790'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("E19").Value = '  +1.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.55'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -0.53%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.152.55'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +1.04%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.412'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -2.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.65'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +1.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.959'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1452'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -0.35%  '
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.958'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -1.09%  '
$ws.Range("E30").Value = '  +1.24%  '
$ws.Range("E31").Value = '  -0.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.402'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -2.54%  '
$ws.Range("E33").Value = '  +0.18%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05231'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +0.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.169'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  -1.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7077'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9996'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -3.74%  '
$ws.Range("E38").Value = '  +0.49%  '
$ws.Range("E39").Value = '  -0.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.728'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +1.86%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9219'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  -2.12%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.140.24'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +6.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4290'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -0.53%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.871'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -2.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.10'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -0.64%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9992'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.34'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -1.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.772'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  +1.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.010.80'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  -0.80%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.186'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +0.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.980'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -1.10%  '
